# Insert a new weekly price record above the current row 72 ("Haba" /
# Vega Modelo de Temuco). This pushes the existing rows 72-97 down to
# 73-98 (dimension grows from A1:R97 to A1:R98), and the new row 72 is
# populated with the latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 72:97 down by one row, opening up a blank row 72.
$ws.Rows(72).Insert()

# Populate the newly inserted row 72 with the new record.
$ws.Range("A72").Value = 10
$ws.Range("B72").Value = "Vega Modelo de Temuco"
$ws.Range("C72").Value = "La Araucanía"
$ws.Range("D72").Value = 45229
$ws.Range("E72").Value = 9
$ws.Range("F72").Value = 100112026
$ws.Range("G72").Value = "Haba"
$ws.Range("H72").Value = "Sin especificar"
$ws.Range("I72").Value = "Primera"
$ws.Range("J72").Value = 265
$ws.Range("K72").Value = 11000
$ws.Range("L72").Value = 13000
$ws.Range("M72").Value = 11491
$ws.Range("N72").Value = "$/saco 25 kilos"
$ws.Range("O72").Value = "Región del Maule"
$ws.Range("P72").Value = 460
$ws.Range("Q72").Value = 25
$ws.Range("R72").Value = "Hortaliza"
